$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the contact's data in row 2 (login/account swapped to "migraineworldsummit")
$ws.Range("D2").Value = "migraineworldsummit"
$ws.Range("F2").Value = "Migraine World Summit"
$ws.Range("K2").Value = 17399
$ws.Range("L2").Value = "do"

# Biography text moves from column O to column P; clear O2.
$ws.Range("P2").Value = "Annual free virtual summit with world leading migraine experts, doctors & specialists.🧠"
$ws.Range("O2").Value = ""

# Drop the old hyperlinks for rows 3/4 (deleting any one hyperlink clears the whole
# collection in this runtime), then re-add the I2/J2 links with the new URLs,
# update the displayed text to match, and restore their Hyperlink cell style.
$ws.Range("I3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("I2"), "https://linktr.ee/Migraineworldsummit")
$ws.Hyperlinks.Add($ws.Range("J2"), "https://www.instagram.com/migraineworldsummit/")
$ws.Range("I2").Value = "https://linktr.ee/Migraineworldsummit"
$ws.Range("J2").Value = "https://www.instagram.com/migraineworldsummit/"
$ws.Range("I2").Style = "Hyperlink"
$ws.Range("J2").Style = "Hyperlink"

# Remove the two extra contact rows entirely.
$ws.Rows("4").Delete()
$ws.Rows("3").Delete()

# Row deletion shrinks the table's bound range automatically; restore it to its
# original extent (A1:P1000) since the table itself isn't meant to change.
$ws.ListObjects.Item(1).Resize($ws.Range("A1:P1000"))
